$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1 -> "placas"
$ws.Range("G1").Value = "placas"

# Row 2
$ws.Range("A2").Value = 1163301040735
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = "GFE450"

# Row 3
$ws.Range("A3").Value = 1633010407032
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = "GFE451"

# Row 4
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = "GFE452"

# Row 5
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "GFE453"

# Row 6
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = "GFE454"

# Row 7
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = "GFE455"

# Row 8
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = "GFE456"

# Row 9
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = "GFE457"

# Row 10
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "GFE458"

# Row 11
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = "GFE459"

# Update selection to A4
$ws.Range("A4").Select() | Out-Null
